# Add modifications to support web element visible and assert data fetch.
# Adds a new "assert_data" worksheet holding the assert keyword/expected-data
# pairs used by the new assertion utility, and refreshes the remembered
# selections on the existing sheets.

$wb = $excel.ActiveWorkbook

# --- customer_data: move remembered selection from H15 to A2 -------------
$customerSheet = $wb.Worksheets.Item("customer_data")
$customerSheet.Activate()
$customerSheet.Range("A2").Select()

# --- search_keyword: move remembered selection from C2 to B2 -------------
$searchSheet = $wb.Worksheets.Item("search_keyword")
$searchSheet.Activate()
$searchSheet.Range("B2").Select()

# --- assert_data: brand new sheet at the end of the workbook -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$assertSheet = $wb.Worksheets.Add($null, $lastSheet)
$assertSheet.Name = "assert_data"

# Cell values are written in this specific order so that new shared-string
# entries land in the same order the source workbook used ("data",
# "search_result", "Results for: ...", "assert_keyword").
$assertSheet.Range("B1").Value = "data"
$assertSheet.Range("A2").Value = "search_result"
$assertSheet.Range("B2").Value = "Results for: EOS Rebel T7 DSLR Camera"
$assertSheet.Range("A1").Value = "assert_keyword"

$assertSheet.Columns.Item(1).ColumnWidth = 20.166666666666664
$assertSheet.Columns.Item(2).ColumnWidth = 74.0

$assertSheet.Range("A1").Select()

# assert_data becomes the active (visible) tab, matching the saved workbook
$assertSheet.Activate()
